# timelog.xlsx update: add four new "direct intersection" work-log
# entries to the "individuals" sheet (rows 151-154). This extends the
# shared duration formula in column D down through row 154 and, after
# recalculation, bumps the "total" sheet's SUM accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("individuals")

# Row 151: 2015-06-26 (serial 42181), 15:30-20:30
$ws.Range("A151").Value = 42181
$ws.Range("B151").Value = 0.64583333333333337
$ws.Range("C151").Value = 0.85416666666666663
$ws.Range("D151").Formula = "=C151-B151"
$ws.Range("E151").Value = "integrated Geometric Tools Library for direct intersection"

# Row 152: 2015-06-29 (serial 42184), 14:00-21:00
$ws.Range("A152").Value = 42184
$ws.Range("B152").Value = 0.58333333333333337
$ws.Range("C152").Value = 0.875
$ws.Range("D152").Formula = "=C152-B152"
$ws.Range("E152").Value = "stabilizing direct intersection"

# Row 153: 2015-06-30 (serial 42185), 09:00-12:00
$ws.Range("A153").Value = 42185
$ws.Range("B153").Value = 0.375
$ws.Range("C153").Value = 0.5
$ws.Range("D153").Formula = "=C153-B153"
$ws.Range("E153").Value = "rewriting direct intersection to operate per cell"

# Row 154: 2015-06-30 (serial 42185), 12:45-18:30
$ws.Range("A154").Value = 42185
$ws.Range("B154").Value = 0.53125
$ws.Range("C154").Value = 0.77083333333333337
$ws.Range("D154").Formula = "=C154-B154"
$ws.Range("E154").Value = "rewriting direct intersection to operate per cell"

# New rows pick up the column-level number formats (date/time/duration)
# automatically, matching A2:D150's existing yyyy-mm-dd / hh:mm / [h]:m
# styles, so no explicit NumberFormat assignment is needed here.

# Move the active selection one row past the new last data row, mirroring
# Excel leaving the cursor just below the used range after data entry.
$ws.Range("E157").Select()

$wb.Save()
